$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 previously held a date value with a date number format (xf idx 1).
# Clear that formatting before repurposing the cell to hold a plain ID string.
$ws.Range("D2").ClearFormats()

# --- Row 1 : headers ----------------------------------------------------
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Localización"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "ID"
$ws.Range("E1").Value = "Tipo"

# --- Row 2 : first record -------------------------------------------------
$ws.Range("A2").Value = "Juan Torres Pardo"
$ws.Range("B2").Value = "40.5N30.99W"
$ws.Range("C2").Value = "juan@example.com"
$ws.Range("D2").Value = "58758L"
$ws.Range("E2").Value = 1

# --- Row 3 : second (new) record ------------------------------------------
$ws.Range("A3").Value = "Juan Torres Pardo"
$ws.Range("B3").Value = "40.5N30.99W"
$ws.Range("C3").Value = "juan@example.com"
$ws.Range("D3").Value = "58758L"
$ws.Range("E3").Value = 1

# The old layout had extra columns (F:I) that are no longer used.
$ws.Range("F1:I2").ClearContents()

# Re-create the mailto hyperlink on C2 (already present) and add a second
# one on C3, each with their own relationship id.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:juan@example.com")

# Column E is no longer the wide "Dirección postal" column - narrow it back.
$ws.Columns("E").ColumnWidth = 10.76

# Match the final selection recorded in the sheet view.
$ws.Range("B12").Select()
